# Actualización automática - registro de operaciones diarias (Semana 2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 24 (Semana 1): nuevo registro diario
$ws.Range("A24").Value = (Get-Date -Year 2025 -Month 8 -Day 19 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C24").Value = 1
$ws.Range("E24").Value = 0

# Fila 28 (Semana 2): nuevo registro diario
$ws.Range("A28").Value = (Get-Date -Year 2025 -Month 8 -Day 22 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = 0

# Fila 29 (Semana 2): nuevo registro diario
$ws.Range("A29").Value = (Get-Date -Year 2025 -Month 8 -Day 26 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C29").Value = 1
$ws.Range("E29").Value = 0

# Fila 30 (Semana 2): nuevo registro diario
$ws.Range("A30").Value = (Get-Date -Year 2025 -Month 8 -Day 28 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C30").Value = 1
$ws.Range("E30").Value = 0

# Recalcular totales y subtotales dependientes de las filas anteriores
$excel.Calculate()

# Dejar la selección activa donde el usuario terminó de capturar datos
$ws.Range("E30:F30").Select()
